# "Updated to Week 6" - refresh contestant statuses + model predictions.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column (N): contestants eliminated this week ---
# shared string 117 "Still in it to Win it" -> 119 "Eliminated"
$ws.Range("N2").Value  = "Eliminated"
$ws.Range("N4").Value  = "Eliminated"
$ws.Range("N9").Value  = "Eliminated"
$ws.Range("N16").Value = "Eliminated"
$ws.Range("N18").Value = "Eliminated"
$ws.Range("N19").Value = "Eliminated"
$ws.Range("N27").Value = "Eliminated"
$ws.Range("N29").Value = "Eliminated"
$ws.Range("N31").Value = "Eliminated"

# --- Linear Regression Prediction column (P) ---
$ws.Range("P8").Value  = 0.177623054172
$ws.Range("P10").Value = 0.177623054172
$ws.Range("P20").Value = 0.177623054172
$ws.Range("P25").Value = 0.13285834223
$ws.Range("P26").Value = 0.177623054172

# --- Random Forest Regression Prediction column (Q) ---
$ws.Range("Q2").Value  = 0.01
$ws.Range("Q5").Value  = 0
$ws.Range("Q7").Value  = 0.162678571429
$ws.Range("Q8").Value  = 0.335
$ws.Range("Q9").Value  = 0.125
$ws.Range("Q10").Value = 0.42
$ws.Range("Q12").Value = 0
$ws.Range("Q14").Value = 0.01
$ws.Range("Q15").Value = 0.01
$ws.Range("Q16").Value = 0
$ws.Range("Q18").Value = 0
$ws.Range("Q20").Value = 0.065
$ws.Range("Q22").Value = 0
$ws.Range("Q24").Value = 0.06
$ws.Range("Q25").Value = 0.21
$ws.Range("Q26").Value = 0.01
$ws.Range("Q29").Value = 0.02
$ws.Range("Q30").Value = 0.0745238095238

# --- Decision Tree Regression Prediction column (R) ---
$ws.Range("R8").Value  = 1
$ws.Range("R26").Value = 0

# --- View state: scroll/selection moved ---
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("T1").Select()
